$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1:J1").EntireColumn.Insert()
$ws.Range("G1").Value = "Informações Adicionais"
$ws.Range("H1").Value = "Resultado (Css)"
$ws.Range("I1").Value = "Resultado (Cl.)"
$ws.Range("J1").Value = "Comentários"
$ws.Range("K1").Value = "Meta"
$ws.Range("H1:J1").EntireColumn.ColumnWidth = $ws.Range("F1").EntireColumn.ColumnWidth()
$ws.Range("G1").EntireColumn.AutoFit()
$w = $ws.Range("L1").EntireColumn.ColumnWidth()
$ws.Range("L1").EntireColumn.ColumnWidth = $w
